$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Debug Nodejs"
$ws.Range("B2").Value = "4h"
$ws.Range("C2").Value = "Bowen"

# "99%" should remain the literal text "99%" (not be auto-converted into the
# number 0.99 with a percentage format), so force the cell to text first.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "99%"
